# Update the cryptos price table (columns B-E) for rows 2-51 to match the
# latest scrape. Price values in column D are stored as literal text in the
# workbook (e.g. "75.926.46", using "." as a thousands separator), so any
# value that COM/Excel would otherwise auto-convert to a number is written
# with a leading apostrophe to force text, then the style is reset back to
# "Normal" so no stray number-format/quote-prefix style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '75.926.46'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '

$ws.Cells.Item(3, 4).Value = '2.886.83'
$ws.Cells.Item(3, 5).Value = '  +6.16%  '

$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).Value = '''195.19'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +4.15%  '

$ws.Cells.Item(6, 4).Value = '''597.44'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +1.03%  '

$ws.Cells.Item(7, 5).Value = '  -0.02%  '

$ws.Cells.Item(8, 4).Value = '''0.554'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +2.06%  '

$ws.Cells.Item(9, 5).Value = '  -2.41%  '

$ws.Cells.Item(10, 4).Value = '2.885.60'
$ws.Cells.Item(10, 5).Value = '  +6.20%  '

$ws.Cells.Item(11, 4).Value = '''0.398'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +9.74%  '

$ws.Cells.Item(12, 5).Value = '  -1.40%  '

$ws.Cells.Item(13, 4).Value = '''4.90'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +2.26%  '

$ws.Cells.Item(14, 4).Value = '3.351.78'
$ws.Cells.Item(14, 5).Value = '  +3.88%  '

$ws.Cells.Item(15, 4).Value = '75.799.42'
$ws.Cells.Item(15, 5).Value = '  +0.34%  '

$ws.Cells.Item(16, 4).Value = '''0.0000190'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.48%  '

$ws.Cells.Item(17, 4).Value = '''27.30'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.78%  '

$ws.Cells.Item(18, 4).Value = '2.878.29'
$ws.Cells.Item(18, 5).Value = '  +5.80%  '

$ws.Cells.Item(19, 4).Value = '''8.89'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -5.64%  '

$ws.Cells.Item(20, 4).Value = '''12.56'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +3.08%  '

$ws.Cells.Item(21, 4).Value = '''377.03'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.60%  '

$ws.Cells.Item(22, 5).Value = '  -0.71%  '

$ws.Cells.Item(23, 4).Value = '''4.15'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.72%  '

$ws.Cells.Item(24, 4).Value = '''71.40'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.32%  '

$ws.Cells.Item(25, 5).Value = '  +0.02%  '

$ws.Cells.Item(26, 4).Value = '3.026.76'
$ws.Cells.Item(26, 5).Value = '  +5.79%  '

$ws.Cells.Item(27, 5).Value = '  -0.41%  '

$ws.Cells.Item(28, 4).Value = '''9.79'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.66%  '

$ws.Cells.Item(29, 5).Value = '  +7.66%  '

$ws.Cells.Item(30, 4).Value = '''0.997'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.34%  '

$ws.Cells.Item(31, 4).Value = '''1.41'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.44%  '

$ws.Cells.Item(32, 4).Value = '''509.27'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -2.79%  '

$ws.Cells.Item(33, 4).Value = '''7.78'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.53%  '

$ws.Cells.Item(34, 4).Value = '''1.81'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.36%  '

$ws.Cells.Item(35, 4).Value = '''0.999'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.02%  '

$ws.Cells.Item(36, 4).Value = '''20.21'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +3.00%  '

$ws.Cells.Item(37, 4).Value = '''163.22'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.82%  '

$ws.Cells.Item(38, 5).Value = '  +1.63%  '

$ws.Cells.Item(39, 5).Value = '  -5.23%  '

$ws.Cells.Item(40, 2).Value = 'Aave'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(40, 4).Value = '''182.24'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +4.24%  '

$ws.Cells.Item(41, 2).Value = 'USDe'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(41, 4).Value = '''1.00'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.00%  '

$ws.Cells.Item(42, 5).Value = '  +3.13%  '

$ws.Cells.Item(43, 4).Value = '''5.00'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.56%  '

$ws.Cells.Item(44, 5).Value = '  -2.48%  '

$ws.Cells.Item(45, 4).Value = '''0.0919'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +7.11%  '

$ws.Cells.Item(46, 5).Value = '  -0.71%  '

$ws.Cells.Item(47, 4).Value = '''40.40'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +3.06%  '

$ws.Cells.Item(48, 4).Value = '''2.36'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -2.27%  '

$ws.Cells.Item(49, 4).Value = '''0.579'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +5.17%  '

$ws.Cells.Item(50, 2).Value = 'Filecoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(50, 4).Value = '''3.76'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.83%  '

$ws.Cells.Item(51, 2).Value = 'Mantle'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51, 4).Value = '''0.664'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +11.04%  '
